# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D) for the eb0053da... file row
# (row 5) on both the zh-cn and de-de worksheets, reflecting the new handoff
# report times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-26 05:32:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-26 05:32:59"
